$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 164
$ws.Cells.Item(2, 9).Value = 162.55556
$ws.Cells.Item(2, 11).Value = 162.55556
$ws.Cells.Item(2, 13).Value = -49.55556000000001

$ws.Cells.Item(41, 8).Value = 1944.45
$ws.Cells.Item(41, 9).Value = 1364.9166
$ws.Cells.Item(41, 10).Value = 2813.75
$ws.Cells.Item(41, 11).Value = 1364.9166
$ws.Cells.Item(41, 12).Value = 2813.75
$ws.Cells.Item(41, 13).Value = -924.9166
$ws.Cells.Item(41, 14).Value = -3693.75

$ws.Cells.Item(62, 8).Value = 5781.222
$ws.Cells.Item(62, 9).Value = 3769.8823
$ws.Cells.Item(62, 11).Value = 3769.8823
$ws.Cells.Item(62, 13).Value = -3145.8823

$ws.Cells.Item(65, 8).Value = 5781.222
$ws.Cells.Item(65, 9).Value = 3769.8823
$ws.Cells.Item(65, 11).Value = 18849.4115
$ws.Cells.Item(65, 13).Value = -15729.4115

$ws.Cells.Item(112, 8).Value = 1616.1578
$ws.Cells.Item(112, 10).Value = 1664.6666
$ws.Cells.Item(112, 12).Value = 4993.9998
$ws.Cells.Item(112, 14).Value = -7209.9998

$ws.Cells.Item(137, 8).Value = 7940.4
$ws.Cells.Item(137, 9).Value = 6249.5
$ws.Cells.Item(137, 11).Value = 18748.5
$ws.Cells.Item(137, 13).Value = -16198.5

$ws.Cells.Item(138, 8).Value = 3668.087
$ws.Cells.Item(138, 9).Value = 3511
$ws.Cells.Item(138, 10).Value = 3729.9697
$ws.Cells.Item(138, 11).Value = 10533
$ws.Cells.Item(138, 12).Value = 11189.9091
$ws.Cells.Item(138, 13).Value = -5393
$ws.Cells.Item(138, 14).Value = -21469.9091

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 6806056.5
$ws.Cells.Item(74, 9).Value = 9262957
$ws.Cells.Item(74, 10).Value = 2332
$ws.Cells.Item(74, 11).Value = 9262957
$ws.Cells.Item(74, 12).Value = 2332
$ws.Cells.Item(74, 13).Value = -9262083
$ws.Cells.Item(74, 14).Value = -4080

$ws.Cells.Item(77, 8).Value = 6806056.5
$ws.Cells.Item(77, 9).Value = 9262957
$ws.Cells.Item(77, 10).Value = 2332
$ws.Cells.Item(77, 11).Value = 46314785
$ws.Cells.Item(77, 12).Value = 11660
$ws.Cells.Item(77, 13).Value = -46310417
$ws.Cells.Item(77, 14).Value = -20396

$ws.Cells.Item(122, 8).Value = 3673.818
$ws.Cells.Item(122, 9).Value = 1482.4
$ws.Cells.Item(122, 11).Value = 4447.200000000001
$ws.Cells.Item(122, 13).Value = -1997.200000000001

$ws.Cells.Item(132, 8).Value = 3951.625
$ws.Cells.Item(132, 9).Value = 1739.8
$ws.Cells.Item(132, 11).Value = 5219.4
$ws.Cells.Item(132, 13).Value = -2689.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 52875.78
$ws.Cells.Item(31, 9).Value = 4459.4287
$ws.Cells.Item(31, 11).Value = 4459.4287
$ws.Cells.Item(31, 13).Value = -4164.4287

$ws.Cells.Item(34, 8).Value = 52875.78
$ws.Cells.Item(34, 9).Value = 4459.4287
$ws.Cells.Item(34, 11).Value = 4459.4287
$ws.Cells.Item(34, 13).Value = -4257.4287

$ws.Cells.Item(58, 8).Value = 6381
$ws.Cells.Item(58, 9).Value = 4445.923
$ws.Cells.Item(58, 10).Value = 9974.714
$ws.Cells.Item(58, 11).Value = 4445.923
$ws.Cells.Item(58, 12).Value = 9974.714
$ws.Cells.Item(58, 13).Value = -4242.923
$ws.Cells.Item(58, 14).Value = -10380.714

$ws.Cells.Item(132, 8).Value = 3500.3513
$ws.Cells.Item(132, 9).Value = 2828.4062
$ws.Cells.Item(132, 10).Value = 7800.8
$ws.Cells.Item(132, 11).Value = 8485.2186
$ws.Cells.Item(132, 12).Value = 23402.4
$ws.Cells.Item(132, 13).Value = -5955.2186
$ws.Cells.Item(132, 14).Value = -28462.4

$ws.Cells.Item(134, 8).Value = 2985.4792
$ws.Cells.Item(134, 9).Value = 2442.2307
$ws.Cells.Item(134, 10).Value = 3627.5
$ws.Cells.Item(134, 11).Value = 7326.6921
$ws.Cells.Item(134, 12).Value = 10882.5
$ws.Cells.Item(134, 13).Value = -4791.6921
$ws.Cells.Item(134, 14).Value = -15952.5

$ws.Cells.Item(136, 8).Value = 6381
$ws.Cells.Item(136, 9).Value = 4445.923
$ws.Cells.Item(136, 10).Value = 9974.714
$ws.Cells.Item(136, 11).Value = 13337.769
$ws.Cells.Item(136, 12).Value = 29924.142
$ws.Cells.Item(136, 13).Value = -10787.769
$ws.Cells.Item(136, 14).Value = -35024.142

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 3023
$ws.Cells.Item(107, 9).Value = 2627.6
$ws.Cells.Item(107, 11).Value = 7882.799999999999
$ws.Cells.Item(107, 13).Value = -5962.799999999999

$ws.Cells.Item(113, 8).Value = 1028.9
$ws.Cells.Item(113, 9).Value = 848.75
$ws.Cells.Item(113, 10).Value = 1749.5
$ws.Cells.Item(113, 11).Value = 2546.25
$ws.Cells.Item(113, 12).Value = 5248.5
$ws.Cells.Item(113, 13).Value = -376.25
$ws.Cells.Item(113, 14).Value = -9588.5

$ws.Cells.Item(114, 8).Value = 1559.25
$ws.Cells.Item(114, 9).Value = 854.8
$ws.Cells.Item(114, 10).Value = 2733.3333
$ws.Cells.Item(114, 11).Value = 2564.4
$ws.Cells.Item(114, 12).Value = 8199.999899999999
$ws.Cells.Item(114, 13).Value = 689.6000000000004
$ws.Cells.Item(114, 14).Value = -14707.9999

$ws.Cells.Item(122, 8).Value = 2071.0625
$ws.Cells.Item(122, 9).Value = 329.5
$ws.Cells.Item(122, 10).Value = 2319.8572
$ws.Cells.Item(122, 11).Value = 2965.5
$ws.Cells.Item(122, 12).Value = 20878.7148
$ws.Cells.Item(122, 13).Value = -515.5
$ws.Cells.Item(122, 14).Value = -25778.7148

$ws.Cells.Item(124, 8).Value = 1212.5
$ws.Cells.Item(124, 9).Value = 1212.5
$ws.Cells.Item(124, 11).Value = 3637.5
$ws.Cells.Item(124, 13).Value = 1272.5

$ws.Cells.Item(128, 8).Value = 203957
$ws.Cells.Item(128, 9).Value = 203957
$ws.Cells.Item(128, 11).Value = 611871
$ws.Cells.Item(128, 13).Value = -606891

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 727.7778
$ws.Cells.Item(97, 10).Value = 946.9
$ws.Cells.Item(97, 12).Value = 946.9
$ws.Cells.Item(97, 14).Value = -1938.9

$ws.Cells.Item(122, 8).Value = 7744.524
$ws.Cells.Item(122, 10).Value = 12599.4
$ws.Cells.Item(122, 12).Value = 37798.2
$ws.Cells.Item(122, 14).Value = -42698.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 5844.3076
$ws.Cells.Item(22, 10).Value = 6699.5
$ws.Cells.Item(22, 12).Value = 6699.5
$ws.Cells.Item(22, 14).Value = -7289.5

$ws.Cells.Item(27, 8).Value = 5844.3076
$ws.Cells.Item(27, 10).Value = 6699.5
$ws.Cells.Item(27, 12).Value = 6699.5
$ws.Cells.Item(27, 14).Value = -6913.5

$ws.Cells.Item(55, 8).Value = 792.44446
$ws.Cells.Item(55, 9).Value = 717.6667
$ws.Cells.Item(55, 10).Value = 942
$ws.Cells.Item(55, 11).Value = 717.6667
$ws.Cells.Item(55, 12).Value = 942
$ws.Cells.Item(55, 13).Value = -544.6667
$ws.Cells.Item(55, 14).Value = -1288

$ws.Cells.Item(61, 8).Value = 1483.4584
$ws.Cells.Item(61, 9).Value = 1480.2
$ws.Cells.Item(61, 10).Value = 1499.75
$ws.Cells.Item(61, 11).Value = 1480.2
$ws.Cells.Item(61, 12).Value = 1499.75
$ws.Cells.Item(61, 13).Value = -1278.2
$ws.Cells.Item(61, 14).Value = -1903.75

$ws.Cells.Item(113, 8).Value = 1483.4584
$ws.Cells.Item(113, 9).Value = 1480.2
$ws.Cells.Item(113, 10).Value = 1499.75
$ws.Cells.Item(113, 11).Value = 1480.2
$ws.Cells.Item(113, 12).Value = 1499.75
$ws.Cells.Item(113, 13).Value = 689.8
$ws.Cells.Item(113, 14).Value = -5839.75

$ws.Cells.Item(136, 8).Value = 7587.294
$ws.Cells.Item(136, 9).Value = 4441.778
$ws.Cells.Item(136, 10).Value = 11126
$ws.Cells.Item(136, 11).Value = 13325.334
$ws.Cells.Item(136, 12).Value = 33378
$ws.Cells.Item(136, 13).Value = -10775.334
$ws.Cells.Item(136, 14).Value = -38478

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1999.625
$ws.Cells.Item(100, 9).Value = 1644.4546
$ws.Cells.Item(100, 11).Value = 3288.9092
$ws.Cells.Item(100, 13).Value = -2747.9092

$ws.Cells.Item(122, 8).Value = 2670.1482
$ws.Cells.Item(122, 9).Value = 2116.5
$ws.Cells.Item(122, 11).Value = 6349.5
$ws.Cells.Item(122, 13).Value = -3899.5

$ws.Cells.Item(136, 8).Value = 6373.4165
$ws.Cells.Item(136, 9).Value = 3366.2632
$ws.Cells.Item(136, 11).Value = 10098.7896
$ws.Cells.Item(136, 13).Value = -7548.7896
